# This workbook holds a weekly time-series of Coliflor prices reported by
# "Macroferia Regional de Talca". A new, more recent weekly record is
# inserted right before the current row 442, pushing all subsequent rows
# (442-518) down by one (to 443-519), exactly as Excel's native row
# insertion does (formats/styles shift along with the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 442; rows 442:518 shift down to 443:519
# and the workbook's used range grows from A1:R518 to A1:R519.
$ws.Rows.Item(442).Insert()

# Populate the newly inserted row 442 with the new weekly record.
$ws.Cells.Item(442, 1).Value = 5
$ws.Cells.Item(442, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(442, 3).Value = "Maule"
$ws.Cells.Item(442, 4).Value = 45218
$ws.Cells.Item(442, 5).Value = 7
$ws.Cells.Item(442, 6).Value = 100112008
$ws.Cells.Item(442, 7).Value = "Coliflor"
$ws.Cells.Item(442, 8).Value = "Sin especificar"
$ws.Cells.Item(442, 9).Value = "Primera"
$ws.Cells.Item(442, 10).Value = 3000
$ws.Cells.Item(442, 11).Value = 800
$ws.Cells.Item(442, 12).Value = 800
$ws.Cells.Item(442, 13).Value = 800
$ws.Cells.Item(442, 14).Value = "`$/unidad"
$ws.Cells.Item(442, 15).Value = "Región del Maule"
$ws.Cells.Item(442, 16).Value = 800
$ws.Cells.Item(442, 17).Value = 1
$ws.Cells.Item(442, 18).Value = "Hortaliza"
